$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 419.66666
$ws.Range("I2").Value = 565
$ws.Range("J2").Value = 129
$ws.Range("K2").Value = 565
$ws.Range("L2").Value = 129
$ws.Range("M2").Value = -452
$ws.Range("N2").Value = -355
$ws.Range("H33").Value = 435.8889
$ws.Range("I33").Value = 367.6
$ws.Range("K33").Value = 367.6
$ws.Range("M33").Value = -138.6
$ws.Range("H70").Value = 16333
$ws.Range("I70").Value = 3999.5
$ws.Range("J70").Value = 22499.75
$ws.Range("K70").Value = 11998.5
$ws.Range("L70").Value = 67499.25
$ws.Range("M70").Value = -11728.5
$ws.Range("N70").Value = -68039.25
$ws.Range("H73").Value = 16333
$ws.Range("I73").Value = 3999.5
$ws.Range("J73").Value = 22499.75
$ws.Range("K73").Value = 11998.5
$ws.Range("L73").Value = 67499.25
$ws.Range("M73").Value = -11062.5
$ws.Range("N73").Value = -69371.25
$ws.Range("H98").Value = 2080.125
$ws.Range("I98").Value = 1971.4667
$ws.Range("J98").Value = 3710
$ws.Range("K98").Value = 1971.4667
$ws.Range("L98").Value = 3710
$ws.Range("M98").Value = -473.4666999999999
$ws.Range("N98").Value = -6706
$ws.Range("H122").Value = 2080.125
$ws.Range("I122").Value = 1971.4667
$ws.Range("J122").Value = 3710
$ws.Range("K122").Value = 5914.4001
$ws.Range("L122").Value = 11130
$ws.Range("M122").Value = -3464.4001
$ws.Range("N122").Value = -16030
$ws.Range("H137").Value = 1916.0714
$ws.Range("I137").Value = 1636
$ws.Range("J137").Value = 2196.1428
$ws.Range("K137").Value = 4908
$ws.Range("L137").Value = 6588.428400000001
$ws.Range("M137").Value = -2358
$ws.Range("N137").Value = -11688.4284
$ws.Range("H138").Value = 3603.3333
$ws.Range("I138").Value = 2988
$ws.Range("J138").Value = 3827.0908
$ws.Range("K138").Value = 8964
$ws.Range("L138").Value = 11481.2724
$ws.Range("M138").Value = -3824
$ws.Range("N138").Value = -21761.2724

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 540953.25
$ws.Range("I2").Value = 1192237.2
$ws.Range("K2").Value = 1192237.2
$ws.Range("M2").Value = -1192124.2
$ws.Range("H45").Value = 1603.125
$ws.Range("I45").Value = 970.1667
$ws.Range("J45").Value = 3502
$ws.Range("K45").Value = 970.1667
$ws.Range("L45").Value = 3502
$ws.Range("M45").Value = -593.1667
$ws.Range("N45").Value = -4256
$ws.Range("H110").Value = 4905.6665
$ws.Range("I110").Value = 5825.3335
$ws.Range("J110").Value = 3986
$ws.Range("K110").Value = 5825.3335
$ws.Range("L110").Value = 3986
$ws.Range("M110").Value = -3780.3335
$ws.Range("N110").Value = -8076
$ws.Range("H116").Value = 540953.25
$ws.Range("I116").Value = 1192237.2
$ws.Range("K116").Value = 1192237.2
$ws.Range("M116").Value = -1189943.2
$ws.Range("H132").Value = 5332.7607
$ws.Range("I132").Value = 4946.7104
$ws.Range("J132").Value = 7166.5
$ws.Range("K132").Value = 14840.1312
$ws.Range("L132").Value = 21499.5
$ws.Range("M132").Value = -12310.1312
$ws.Range("N132").Value = -26559.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 540953.25
$ws.Range("I3").Value = 1192237.2
$ws.Range("K3").Value = 1192237.2
$ws.Range("M3").Value = -1192123.2
$ws.Range("H53").Value = 61166.332
$ws.Range("J53").Value = 61166.332
$ws.Range("L53").Value = 61166.332
$ws.Range("N53").Value = -62314.332
$ws.Range("H100").Value = 33585.89
$ws.Range("J100").Value = 33585.89
$ws.Range("L100").Value = 33585.89
$ws.Range("N100").Value = -35749.89
$ws.Range("H107").Value = 1996.125
$ws.Range("I107").Value = 2012.8334
$ws.Range("J107").Value = 1946
$ws.Range("K107").Value = 2012.8334
$ws.Range("L107").Value = 1946
$ws.Range("M107").Value = -92.83339999999998
$ws.Range("N107").Value = -5786

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13762.909
$ws.Range("J31").Value = 16204.111
$ws.Range("L31").Value = 16204.111
$ws.Range("N31").Value = -16794.111
$ws.Range("H34").Value = 13762.909
$ws.Range("J34").Value = 16204.111
$ws.Range("L34").Value = 16204.111
$ws.Range("N34").Value = -16608.111
$ws.Range("H68").Value = 100740.164
$ws.Range("J68").Value = 96124
$ws.Range("L68").Value = 96124
$ws.Range("N68").Value = -97622
$ws.Range("H71").Value = 100740.164
$ws.Range("J71").Value = 96124
$ws.Range("L71").Value = 288372
$ws.Range("N71").Value = -295860
$ws.Range("H86").Value = 3965
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 3965
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H99").Value = 3415.7222
$ws.Range("I99").Value = 3489.6667
$ws.Range("K99").Value = 3489.6667
$ws.Range("M99").Value = -1991.6667
$ws.Range("H107").Value = 838837.75
$ws.Range("I107").Value = 1553387.6
$ws.Range("J107").Value = 5196.3335
$ws.Range("K107").Value = 1553387.6
$ws.Range("L107").Value = 5196.3335
$ws.Range("M107").Value = -1551467.6
$ws.Range("N107").Value = -9036.333500000001
$ws.Range("H126").Value = 3415.7222
$ws.Range("I126").Value = 3489.6667
$ws.Range("K126").Value = 10469.0001
$ws.Range("M126").Value = -7999.000100000001
$ws.Range("H135").Value = 99996.5
$ws.Range("J135").Value = 99996.5
$ws.Range("L135").Value = 99996.5
$ws.Range("N135").Value = -110136.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 737083
$ws.Range("I4").Value = 778026.5600000001
$ws.Range("J4").Value = 99
$ws.Range("K4").Value = 2334079.68
$ws.Range("L4").Value = 297
$ws.Range("M4").Value = -2333967.68
$ws.Range("N4").Value = -521
$ws.Range("H76").Value = 4500
$ws.Range("I76").Value = 4500
$ws.Range("K76").Value = 13500
$ws.Range("M76").Value = -13117
$ws.Range("H79").Value = 4500
$ws.Range("I79").Value = 4500
$ws.Range("K79").Value = 13500
$ws.Range("M79").Value = -12174
$ws.Range("H92").Value = 395
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H113").Value = 53797.844
$ws.Range("I113").Value = 112192.336
$ws.Range("K113").Value = 336577.008
$ws.Range("M113").Value = -334407.008

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 55749.25
$ws.Range("J105").Value = 55999.668
$ws.Range("L105").Value = 55999.668
$ws.Range("N105").Value = -62987.668
$ws.Range("H113").Value = 3110.3704
$ws.Range("I113").Value = 2262.8096
$ws.Range("K113").Value = 2262.8096
$ws.Range("M113").Value = -92.80960000000005
$ws.Range("H122").Value = 64496.42
$ws.Range("I122").Value = 81309.664
$ws.Range("J122").Value = 1446.75
$ws.Range("K122").Value = 243928.992
$ws.Range("L122").Value = 4340.25
$ws.Range("M122").Value = -241478.992
$ws.Range("N122").Value = -9240.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6226.6206
$ws.Range("I7").Value = 5010.273
$ws.Range("K7").Value = 5010.273
$ws.Range("M7").Value = -4898.273
$ws.Range("H40").Value = 8584.75
$ws.Range("I40").Value = 8298.117
$ws.Range("K40").Value = 8298.117
$ws.Range("M40").Value = -8162.117
$ws.Range("H46").Value = 1221.3334
$ws.Range("I46").Value = 999
$ws.Range("J46").Value = 3000
$ws.Range("K46").Value = 999
$ws.Range("L46").Value = 3000
$ws.Range("M46").Value = -811
$ws.Range("N46").Value = -3376
$ws.Range("H68").Value = 2317035.5
$ws.Range("I68").Value = 2606040
$ws.Range("J68").Value = 5000
$ws.Range("K68").Value = 2606040
$ws.Range("L68").Value = 5000
$ws.Range("M68").Value = -2605291
$ws.Range("N68").Value = -6498
$ws.Range("H71").Value = 2317035.5
$ws.Range("I71").Value = 2606040
$ws.Range("J71").Value = 5000
$ws.Range("K71").Value = 13030200
$ws.Range("L71").Value = 25000
$ws.Range("M71").Value = -13026456
$ws.Range("N71").Value = -32488
$ws.Range("H105").Value = 29999.5
$ws.Range("J105").Value = 29999.5
$ws.Range("L105").Value = 29999.5
$ws.Range("N105").Value = -36987.5
$ws.Range("H126").Value = 6226.6206
$ws.Range("I126").Value = 5010.273
$ws.Range("K126").Value = 15030.819
$ws.Range("M126").Value = -12560.819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 25750
$ws.Range("I62").Value = 21000
$ws.Range("K62").Value = 21000
$ws.Range("M62").Value = -20376
$ws.Range("H65").Value = 25750
$ws.Range("I65").Value = 21000
$ws.Range("K65").Value = 105000
$ws.Range("M65").Value = -101880
$ws.Range("H126").Value = 1783.2142
$ws.Range("I126").Value = 1913.3334
$ws.Range("J126").Value = 1002.5
$ws.Range("K126").Value = 5740.0002
$ws.Range("L126").Value = 3007.5
$ws.Range("M126").Value = -3270.0002
$ws.Range("N126").Value = -7947.5
$ws.Range("H132").Value = 5739.475
$ws.Range("I132").Value = 5691.8237
$ws.Range("K132").Value = 17075.4711
$ws.Range("M132").Value = -14545.4711
